# The diff shows that old row 14 (Deposit / Credit Card / Sipay / 216.72),
# which was formatted with a red font / red fill highlight, was removed
# entirely. Rows 15-20 shift up to become rows 14-19 (their contents are
# unchanged), the sheet's used dimension shrinks from A1:AB20 to A1:AB19,
# and the now-unreferenced "Sipay" shared string drops out of the saved
# string table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the row's formatting first so the red-highlight style is no longer
# referenced by any cell, then delete the whole row (shifting rows below
# it up by one), matching a standard "right-click row header > Delete".
$ws.Rows.Item(14).ClearFormats() | Out-Null
$ws.Rows.Item(14).Delete() | Out-Null

# Leave the selection where the author's cursor ended up after the edit.
$ws.Range("K26").Select() | Out-Null
